$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B35").Value = 6445.759765625
$ws.Range("C35").Value = 6446.5498046875
$ws.Range("F35").Value = 5135300000
$ws.Range("G35").Value = 0.0113454358466127
